# Insert a new weekly price record as row 50 ("Fruta / hortaliza, semanal").
# Excel's native row-insert shifts the existing rows 50-119 down to 51-120
# (and grows the sheet dimension to A1:R120 automatically), exactly like a
# user selecting row 50 and choosing Insert > Entire Row before typing the
# new record in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Insert()

$ws.Range("A50").Value = 9
$ws.Range("B50").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C50").Value = "Metropolitana"
$ws.Range("D50").Value = 44771
$ws.Range("E50").Value = 13
$ws.Range("F50").Value = 100112022
$ws.Range("G50").Value = "Arveja Verde"
$ws.Range("H50").Value = "Perfection"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 25
$ws.Range("K50").Value = 43000
$ws.Range("L50").Value = 43000
$ws.Range("M50").Value = 43000
$ws.Range("N50").Value = '$/malla 25 kilos'
$ws.Range("O50").Value = "Provincia de Huasco"
$ws.Range("P50").Value = 1720
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
